$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.008007
$ws.Range("H2").Value = 0.024021
$ws.Range("M2").Value = 0.8562376666666668
$ws.Range("N2").Value = 2.568713
$ws.Range("O2").Value = 0.09910590237923185
$ws.Range("P2").Value = 0.09910590237923182
$ws.Range("Q2").Value = 0.006855894997000001
$ws.Range("R2").Value = 0.06170305497300001
$ws.Range("S2").Value = 0.09910590237923185
$ws.Range("T2").Value = 0.09910590237923182
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.008007
$ws.Range("H3").Value = 0.024021
$ws.Range("O3").Value = 0.448651387965602
$ws.Range("P3").Value = 0.4486513879656019
$ws.Range("Q3").Value = 0.031036565253
$ws.Range("R3").Value = 0.279329087277
$ws.Range("S3").Value = 0.448651387965602
$ws.Range("T3").Value = 0.4486513879656019
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.008007
$ws.Range("H4").Value = 0.024021
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.0338
$ws.Range("N4").Value = 0.1014
$ws.Range("O4").Value = 0.003912207592383465
$ws.Range("P4").Value = 0.003912207592383464
$ws.Range("Q4").Value = 0.0002706366
$ws.Range("R4").Value = 0.0024357294
$ws.Range("S4").Value = 0.003912207592383465
$ws.Range("T4").Value = 0.003912207592383464
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.008007
$ws.Range("H5").Value = 0.024021
$ws.Range("M5").Value = 3.841810666666666
$ws.Range("N5").Value = 11.525432
$ws.Range("O5").Value = 0.4446733981844117
$ws.Range("P5").Value = 0.4446733981844115
$ws.Range("Q5").Value = 0.030761378008
$ws.Range("R5").Value = 0.276852402072
$ws.Range("S5").Value = 0.4446733981844117
$ws.Range("T5").Value = 0.4446733981844115
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.008007
$ws.Range("H6").Value = 0.024021
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.01560133333333333
$ws.Range("N6").Value = 0.046804
$ws.Range("O6").Value = 0.001805788601123429
$ws.Range("P6").Value = 0.001805788601123429
$ws.Range("Q6").Value = 0.000124919876
$ws.Range("R6").Value = 0.001124278884
$ws.Range("S6").Value = 0.001805788601123429
$ws.Range("T6").Value = 0.001805788601123429
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.008007
$ws.Range("H7").Value = 0.024021
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.01599466666666667
$ws.Range("N7").Value = 0.047984
$ws.Range("O7").Value = 0.001851315277247813
$ws.Range("P7").Value = 0.001851315277247812
$ws.Range("Q7").Value = 0.000128069296
$ws.Range("R7").Value = 0.001152623664
$ws.Range("S7").Value = 0.001851315277247813
$ws.Range("T7").Value = 0.001851315277247812
